$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header labels for columns H and I
$ws.Range("H1").Value = "HW_BW_In"
$ws.Range("I1").Value = "KW_BW_In"

# New data values for rows 2-13, columns H (HW_BW_In) and I (KW_BW_In)
$values = @(
    @(5.0999999999999996, 12.3),
    @(4.3, 11.5),
    @(4.3, 11),
    @(4.9000000000000004, 11.6),
    @(4.2, 11.1),
    @(4.3, 12),
    @(5.4, 12.9),
    @(4.2, 11.3),
    @(3, 9.8000000000000007),
    @(2.8, 9),
    @(3.3, 10.8),
    @(3.1, 10.1)
)

$row = 2
foreach ($pair in $values) {
    $ws.Cells.Item($row, 8).Value = $pair[0]
    $ws.Cells.Item($row, 9).Value = $pair[1]
    $row++
}

# Update the selected cell to match the final state (I14)
$ws.Range("I14").Select()
